$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap B24/D24 with B25/D25 (A340V/6 <-> A340T/4)
$b24 = $ws.Range("B24").Value2
$d24 = $ws.Range("D24").Value2
$b25 = $ws.Range("B25").Value2
$d25 = $ws.Range("D25").Value2

$ws.Range("B24").Value2 = $b25
$ws.Range("D24").Value2 = $d25
$ws.Range("B25").Value2 = $b24
$ws.Range("D25").Value2 = $d24

# Swap B28/D28 with B29/D29 (L390Q/6 <-> L390M/3)
$b28 = $ws.Range("B28").Value2
$d28 = $ws.Range("D28").Value2
$b29 = $ws.Range("B29").Value2
$d29 = $ws.Range("D29").Value2

$ws.Range("B28").Value2 = $b29
$ws.Range("D28").Value2 = $d29
$ws.Range("B29").Value2 = $b28
$ws.Range("D29").Value2 = $d28

# Adjust view: move selection to J12 (also clears the scrolled topLeftCell)
$ws.Activate()
$ws.Range("J12").Select()
